$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.8
$ws.Range("G2").Value = 3.25
$ws.Range("L2").Value = 1.01
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 3.5
$ws.Range("O2").Value = 1.32
$ws.Range("R2").Value = 1.33
$ws.Range("S2").Value = 3.3
$ws.Range("T2").Value = 1.62
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 1.57
$ws.Range("W2").Value = 1.44
$ws.Range("X2").Value = 17
$ws.Range("Y2").Value = 11.5
$ws.Range("Z2").Value = 18
$ws.Range("AA2").Value = 40
$ws.Range("AB2").Value = 12.5
$ws.Range("AC2").Value = 8.4
$ws.Range("AD2").Value = 12.5
$ws.Range("AE2").Value = 30
$ws.Range("AF2").Value = 22
$ws.Range("AG2").Value = 14
$ws.Range("AH2").Value = 18
$ws.Range("AI2").Value = 44
$ws.Range("AJ2").Value = 55
$ws.Range("AK2").Value = 36
$ws.Range("AL2").Value = 48
$ws.Range("AM2").Value = 110
$ws.Range("AN2").Value = 32
$ws.Range("AO2").Value = 25

# Row 3
$ws.Range("N3").Value = 5.6
$ws.Range("P3").Value = 2.48
$ws.Range("S3").Value = 2.54
$ws.Range("Y3").Value = 24

# Row 4
$ws.Range("H4").Value = 1.33
$ws.Range("J4").Value = 5.4
$ws.Range("S4").Value = 2.52

# Row 5
$ws.Range("I5").Value = 5.4

# Row 6
$ws.Range("G6").Value = 1.39
$ws.Range("H6").Value = 9.4
$ws.Range("K6").Value = 6
$ws.Range("R6").Value = 1.76
$ws.Range("T6").Value = 1.75

# Row 7
$ws.Range("H7").Value = 2.1
$ws.Range("K7").Value = 3.9
$ws.Range("N7").Value = 5.2

# Row 8
$ws.Range("F8").Value = 1.71
$ws.Range("P8").Value = 1.9
$ws.Range("AJ8").Value = 17
$ws.Range("AL8").Value = 42
